# Auto commit at 2025-09-13  8:40:18.74
# Append the two new daily data rows (row 86: 四方坪站, row 87: 高岭站) for
# date 2025-09-12 (serial 45912) to the bottom of the "Sheet1" data table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 86 : 四方坪站 -------------------------------------------------
$ws.Cells.Item(86, 1).Value = 45912
$ws.Cells.Item(86, 2).Value = "四方坪站"
$ws.Cells.Item(86, 3).Value = 11336.09
$ws.Cells.Item(86, 4).Value = 9314.5499999999993
$ws.Cells.Item(86, 5).Value = 3898.65
$ws.Cells.Item(86, 6).Value = 464

# --- Row 87 : 高岭站 ----------------------------------------------------
$ws.Cells.Item(87, 1).Value = 45912
$ws.Cells.Item(87, 2).Value = "高岭站"
$ws.Cells.Item(87, 3).Value = 5380.24
$ws.Cells.Item(87, 4).Value = 4293.22
$ws.Cells.Item(87, 5).Value = 1348.35
$ws.Cells.Item(87, 6).Value = 178

# Re-use the exact cell styles from the row above (date style in column A,
# integer style in column F) instead of re-typing a NumberFormat string, so
# the existing style indices in styles.xml are reused rather than cloned.
$ws.Range("A85:F85").Copy()
$ws.Range("A86:F86").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("A87:F87").PasteSpecial(-4122)   # xlPasteFormats

# Move the view down to the new bottom of the sheet and put the selection
# on H87, matching where the user's cursor ended up after typing the data.
$win = $excel.ActiveWindow
$win.ScrollRow = 76
$win.ScrollColumn = 1
$ws.Range("H87").Select()

$wb.Save()
